$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 3 values
$ws.Range("A3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

# Update row 4 values
$ws.Range("B4").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# Update the selected/active cell to B1
$ws.Activate()
$ws.Range("B1").Select()
